$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data block (rows 818..930) down by 3 rows (to 821..933)
# to make room for 3 new rows of data inserted at 818..820.
$src = $ws.Range("A818:R930")
$vals = $src.Value()
$dst = $ws.Range("A821:R933")
$dst.Value = $vals

# Write the 3 new rows (date 2022-06-06 / serial 44748, "Cuatro cascos verde")
# at the now-vacated 818..820 slot.
$ws.Range("D818").Value = 44748
$ws.Range("H818").Value = "Cuatro cascos verde"
$ws.Range("J818").Value = 1100
$ws.Range("K818").Value = 29000
$ws.Range("L818").Value = 30000
$ws.Range("M818").Value = 29500
$ws.Range("P818").Value = 1639

$ws.Range("D819").Value = 44748
$ws.Range("H819").Value = "Cuatro cascos verde"
$ws.Range("J819").Value = 700
$ws.Range("K819").Value = 26000
$ws.Range("L819").Value = 27000
$ws.Range("M819").Value = 26500
$ws.Range("P819").Value = 1472

$ws.Range("D820").Value = 44748
$ws.Range("H820").Value = "Cuatro cascos verde"
$ws.Range("J820").Value = 500
$ws.Range("K820").Value = 23000
$ws.Range("L820").Value = 24000
$ws.Range("M820").Value = 23500
$ws.Range("P820").Value = 1306
